# credenciales.xlsx — replace outgoing teammate's credential row with the
# new owner, and lay out a blank 10x3 block (rows 5-14, cols A-C) that the
# download-report robot/main script will populate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("credenciales")

# The old hyperlinked mailto address in A2 is replaced by a plain-text
# e-mail for the new user (natalia.gonzalezb@griky.co); drop the hyperlink
# entirely since the new value is just text, not a mailto link.
[void]$ws.Hyperlinks.Delete()
$ws.Range("A2").Value = "natalia.gonzalezb@griky.co"

# Match the plain, non-link look (the old "Hipervinculo" font was blue &
# underlined; the replacement text uses a plain small grey font).
$f = $ws.Range("A2").Font
$f.Name = "Segoe UI"
$f.Size = 8
$f.Color = 2368548   # RGB(36,36,36) -> FF242424
$f.Underline = -4142 # xlUnderlineStyleNone

# Reserve a blank 10-row x 3-col area (A5:C14) for the new report output.
[void]($ws.Range("A5:C14").Font.Name = "Calibri")

[void]$ws.Range("D4").Select()
